# Oppdater felter-faner i test-kodebok
#
# "1-Testskjema-felter" and "2-Sluttskjema-felter" were still empty
# placeholder sheets; this fills both with the same field-header block
# (A1:F14) that already lives on "1-Testskjema", and then updates the
# active-sheet/selection bookkeeping so "2-Sluttskjema-felter" is the
# sheet that ends up selected (matching the authored workbook state).

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("1-Testskjema")
$srcRange = $template.Range("A1:F14")

$targets = @(
    $wb.Worksheets.Item("1-Testskjema-felter"),
    $wb.Worksheets.Item("2-Sluttskjema-felter")
)

foreach ($dst in $targets) {
    # Formats first (styles/number formats/borders), then values/strings,
    # so the destination ends up with the same cell styles (s="5" header
    # row, s="4" data rows, s="1"/s="2" blank footer row) as the source.
    $srcRange.Copy()
    [void]$dst.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
    $srcRange.Copy()
    [void]$dst.Range("A1").PasteSpecial(-4163)   # xlPasteValues
    [void]$dst.Range("A1:F14").Select()
}

$excel.CutCopyMode = 0

# The template sheet's selection moves off the old D21 cell and onto the
# block that was just copied from it.
[void]$template.Range("A1:F14").Select()

# "2-Sluttskjema-felter" becomes the active / selected tab.
$felter2 = $targets[1]
$felter2.Activate()
[void]$felter2.Range("A1:F14").Select()
